# Trap card data update:
#  - Strain-after-movement effect (Delayed bomb) is no longer a separate rule;
#    folded into the card's own effect text.
#  - Most trap cards that used to only affect "this card's own slot" now also
#    affect the opposite slot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tomb (墓碑) - D2: "顶端" -> "第1张"
$ws.Range("D2").Value = "回合结束时：将主牌堆顶2张牌送墓。<br>
开战时：用墓地第1张牌替换本牌。"

# Spike trap (尖刺) - D5: also affects opposite slot
$ws.Range("D5").Value = "回合结束时：横置本牌，然后本牌所在槽位和对位槽位中所有其他牌点数减1。"

# Cryogas vent (冷气喷口) - D6: also affects opposite slot
$ws.Range("D6").Value = "回合结束时：横置本牌所在槽位和对位槽位中所有牌。"

# Mine (地雷) - D11: also affects opposite slot
$ws.Range("D11").Value = "有牌移入本牌所在槽位时：本牌所在槽位和对位槽位中所有牌点数减1。"

# Delayed bomb (延迟爆弹) - D13: also affects opposite slot
$ws.Range("D13").Value = "回合结束时：移动到1个相邻槽位，然后点数减1，本牌点数因此降至0时，消灭本牌所在槽位或对位槽位中的1张其他牌。"

# Secret door (暗门) - D14: strain-after-movement rule replaced with explicit effect text
$ws.Range("D14").Value = "回合结束时：同槽位中有怪物牌时，将本牌和同槽位中1张怪物牌移动到对位槽位。"

# Update the remembered selection to match the authored workbook (D15)
$null = $ws.Range("D15").Select()
